$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3366.2856
$ws.Range("I64").Value = 3343.75
$ws.Range("K64").Value = 3343.75
$ws.Range("M64").Value = -3095.75
$ws.Range("H67").Value = 3366.2856
$ws.Range("I67").Value = 3343.75
$ws.Range("K67").Value = 3343.75
$ws.Range("M67").Value = -2485.75
$ws.Range("H98").Value = 239021.31
$ws.Range("I98").Value = 761.80554
$ws.Range("J98").Value = 1668578.4
$ws.Range("K98").Value = 761.80554
$ws.Range("L98").Value = 1668578.4
$ws.Range("M98").Value = 736.19446
$ws.Range("N98").Value = -1671574.4
$ws.Range("H122").Value = 239021.31
$ws.Range("I122").Value = 761.80554
$ws.Range("J122").Value = 1668578.4
$ws.Range("K122").Value = 2285.41662
$ws.Range("L122").Value = 5005735.199999999
$ws.Range("M122").Value = 164.58338
$ws.Range("N122").Value = -5010635.199999999
$ws.Range("H135").Value = 1612.2307
$ws.Range("I135").Value = 1984.4445
$ws.Range("J135").Value = 774.75
$ws.Range("K135").Value = 17860.0005
$ws.Range("L135").Value = 6972.75
$ws.Range("M135").Value = -15325.0005
$ws.Range("N135").Value = -12042.75
$ws.Range("H137").Value = 60232.895
$ws.Range("I137").Value = 93032.914
$ws.Range("J137").Value = 4004.2856
$ws.Range("K137").Value = 279098.742
$ws.Range("L137").Value = 12012.8568
$ws.Range("M137").Value = -276548.742
$ws.Range("N137").Value = -17112.8568
$ws.Range("H138").Value = 3202.718
$ws.Range("J138").Value = 3541.2334
$ws.Range("L138").Value = 10623.7002
$ws.Range("N138").Value = -20903.7002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4629.4033
$ws.Range("I32").Value = 2405
$ws.Range("J32").Value = 16492.889
$ws.Range("K32").Value = 2405
$ws.Range("L32").Value = 16492.889
$ws.Range("M32").Value = -2118
$ws.Range("N32").Value = -17066.889
$ws.Range("H61").Value = 5519.05
$ws.Range("I61").Value = 3890.1667
$ws.Range("J61").Value = 7962.375
$ws.Range("K61").Value = 3890.1667
$ws.Range("L61").Value = 7962.375
$ws.Range("M61").Value = -3678.1667
$ws.Range("N61").Value = -8386.375
$ws.Range("H136").Value = 5519.05
$ws.Range("I136").Value = 3890.1667
$ws.Range("J136").Value = 7962.375
$ws.Range("K136").Value = 11670.5001
$ws.Range("L136").Value = 23887.125
$ws.Range("M136").Value = -9120.500100000001
$ws.Range("N136").Value = -28987.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 6502.6665
$ws.Range("I16").Value = 6004
$ws.Range("J16").Value = 7500
$ws.Range("K16").Value = 6004
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = -5834
$ws.Range("N16").Value = -7840
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H134").Value = 3271.6611
$ws.Range("I134").Value = 3634.8538
$ws.Range("K134").Value = 10904.5614
$ws.Range("M134").Value = -8369.561399999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 538315.5
$ws.Range("I99").Value = 1434999.8
$ws.Range("K99").Value = 1434999.8
$ws.Range("M99").Value = -1433501.8
$ws.Range("H126").Value = 538315.5
$ws.Range("I126").Value = 1434999.8
$ws.Range("K126").Value = 4304999.4
$ws.Range("M126").Value = -4302529.4
$ws.Range("H134").Value = 3688.2632
$ws.Range("I134").Value = 3372.5625
$ws.Range("J134").Value = 5372
$ws.Range("K134").Value = 10117.6875
$ws.Range("L134").Value = 16116
$ws.Range("M134").Value = -7582.6875
$ws.Range("N134").Value = -21186

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 6500
$ws.Range("I11").Value = 5000
$ws.Range("J11").Value = 8000
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = -14860
$ws.Range("N11").Value = -24280
$ws.Range("H34").Value = 1116.697
$ws.Range("I34").Value = 654.08
$ws.Range("J34").Value = 2562.375
$ws.Range("K34").Value = 1962.24
$ws.Range("L34").Value = 7687.125
$ws.Range("M34").Value = -1878.24
$ws.Range("N34").Value = -7855.125
$ws.Range("H39").Value = 1497.1875
$ws.Range("I39").Value = 505.25
$ws.Range("J39").Value = 2489.125
$ws.Range("K39").Value = 1515.75
$ws.Range("L39").Value = 7467.375
$ws.Range("M39").Value = -1221.75
$ws.Range("N39").Value = -8055.375
$ws.Range("H55").Value = 7044.933
$ws.Range("J55").Value = 7044.933
$ws.Range("L55").Value = 21134.799
$ws.Range("N55").Value = -21488.799
$ws.Range("H75").Value = 2833.3333
$ws.Range("J75").Value = 3500
$ws.Range("L75").Value = 10500
$ws.Range("N75").Value = -12496
$ws.Range("H78").Value = 2833.3333
$ws.Range("J78").Value = 3500
$ws.Range("L78").Value = 31500
$ws.Range("N78").Value = -41484
$ws.Range("H131").Value = 13159568
$ws.Range("I131").Value = 62501216
$ws.Range("K131").Value = 187503648
$ws.Range("M131").Value = -187498608
$ws.Range("H138").Value = 1376.2
$ws.Range("I138").Value = 1376.2
$ws.Range("K138").Value = 4128.6
$ws.Range("M138").Value = 1011.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3336105.2
$ws.Range("I3").Value = 3507.6365
$ws.Range("J3").Value = 12500748
$ws.Range("K3").Value = 3507.6365
$ws.Range("L3").Value = 12500748
$ws.Range("M3").Value = -3391.6365
$ws.Range("N3").Value = -12500980
$ws.Range("H104").Value = 61500
$ws.Range("J104").Value = 61500
$ws.Range("L104").Value = 61500
$ws.Range("N104").Value = -68488
$ws.Range("H126").Value = 9437.888999999999
$ws.Range("I126").Value = 10501.363
$ws.Range("K126").Value = 31504.089
$ws.Range("M126").Value = -29034.089

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2674.2144
$ws.Range("I61").Value = 2286.6667
$ws.Range("K61").Value = 2286.6667
$ws.Range("M61").Value = -2084.6667
$ws.Range("H113").Value = 2674.2144
$ws.Range("I113").Value = 2286.6667
$ws.Range("K113").Value = 2286.6667
$ws.Range("M113").Value = -116.6667000000002
$ws.Range("H122").Value = 9999.5
$ws.Range("I122").Value = 9999.5
$ws.Range("K122").Value = 29998.5
$ws.Range("M122").Value = -27548.5
$ws.Range("H132").Value = 4179.6
$ws.Range("I132").Value = 3926.2104
$ws.Range("K132").Value = 11778.6312
$ws.Range("M132").Value = -9248.6312

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4200.8237
$ws.Range("I122").Value = 4200.8237
$ws.Range("K122").Value = 12602.4711
$ws.Range("M122").Value = -10152.4711
$ws.Range("H136").Value = 241774.22
$ws.Range("I136").Value = 281285.25
$ws.Range("J136").Value = 4708
$ws.Range("K136").Value = 843855.75
$ws.Range("L136").Value = 14124
$ws.Range("M136").Value = -841305.75
$ws.Range("N136").Value = -19224
